$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage (matches source workbook, which stores all of these as inline strings,
# including numeric-looking Price values) by prefixing with a leading apostrophe so Excel
# does not auto-convert numeric-looking strings (prices, percentages) into numbers.

$ws.Range("D2").Value = "'70.872.12"
$ws.Range("E2").Value = "'  +3.26%  "
$ws.Range("D3").Value = "'3.797.70"
$ws.Range("E3").Value = "'  +1.35%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("D5").Value = "'697.17"
$ws.Range("E5").Value = "'  +11.47%  "
$ws.Range("D6").Value = "'173.33"
$ws.Range("E6").Value = "'  +5.31%  "
$ws.Range("D7").Value = "'3.796.06"
$ws.Range("E7").Value = "'  +1.38%  "
$ws.Range("E8").Value = "'  -0.06%  "
$ws.Range("E9").Value = "'  +1.39%  "
$ws.Range("E10").Value = "'  +3.54%  "
$ws.Range("D11").Value = "'7.44"
$ws.Range("E11").Value = "'  +8.09%  "
$ws.Range("E12").Value = "'  +1.52%  "
$ws.Range("E13").Value = "'  +9.40%  "
$ws.Range("D14").Value = "'36.27"
$ws.Range("E14").Value = "'  +4.85%  "
$ws.Range("D15").Value = "'4.438.23"
$ws.Range("E15").Value = "'  +1.35%  "
$ws.Range("D16").Value = "'3.797.71"
$ws.Range("E16").Value = "'  +1.55%  "
$ws.Range("D17").Value = "'70.868.19"
$ws.Range("E17").Value = "'  +3.27%  "
$ws.Range("D18").Value = "'17.81"
$ws.Range("E19").Value = "'  +3.82%  "
$ws.Range("E20").Value = "'  +0.63%  "
$ws.Range("D21").Value = "'11.09"
$ws.Range("E21").Value = "'  +17.69%  "
$ws.Range("D22").Value = "'483.11"
$ws.Range("E22").Value = "'  +3.39%  "
$ws.Range("E23").Value = "'  +2.28%  "
$ws.Range("E24").Value = "'  +3.75%  "
$ws.Range("E25").Value = "'  +2.49%  "
$ws.Range("E26").Value = "'  +2.94%  "
$ws.Range("D27").Value = "'10.44"
$ws.Range("E27").Value = "'  +4.48%  "
$ws.Range("E28").Value = "'  +3.60%  "
$ws.Range("D29").Value = "'3.947.95"
$ws.Range("E29").Value = "'  +1.36%  "
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "'  -0.02%  "
$ws.Range("D31").Value = "'3.05"
$ws.Range("E31").Value = "'  +15.70%  "
$ws.Range("E32").Value = "'  +6.27%  "
$ws.Range("D33").Value = "'2.28"
$ws.Range("E33").Value = "'  +1.36%  "
$ws.Range("D34").Value = "'29.58"
$ws.Range("E34").Value = "'  +4.85%  "
$ws.Range("D35").Value = "'0.181"
$ws.Range("E35").Value = "'  +2.80%  "
$ws.Range("E36").Value = "'  +4.50%  "
$ws.Range("E37").Value = "'  +0.15%  "
$ws.Range("D38").Value = "'3.747.91"
$ws.Range("E38").Value = "'  +1.31%  "
$ws.Range("E39").Value = "'  +2.93%  "
$ws.Range("D40").Value = "'3.50"
$ws.Range("E40").Value = "'  +9.14%  "
$ws.Range("E41").Value = "'  +4.70%  "
$ws.Range("D42").Value = "'2.22"
$ws.Range("E42").Value = "'  +14.39%  "
$ws.Range("B43").Value = "'FLOKI"
$ws.Range("C43").Value = "'https://coinranking.com/coin/fmHk13Rqw+floki-floki"
$ws.Range("D43").Value = "'0.000327"
$ws.Range("E43").Value = "'  +24.30%  "
$ws.Range("B44").Value = "'Mantle"
$ws.Range("C44").Value = "'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D44").Value = "'0.971"
$ws.Range("E44").Value = "'  +2.17%  "
$ws.Range("D45").Value = "'1.00"
$ws.Range("E45").Value = "'  +0.04%  "
$ws.Range("E46").Value = "'  +0.00%  "
$ws.Range("D47").Value = "'162.82"
$ws.Range("E47").Value = "'  +4.16%  "
$ws.Range("D48").Value = "'49.09"
$ws.Range("E48").Value = "'  +3.88%  "
$ws.Range("D49").Value = "'44.83"
$ws.Range("E49").Value = "'  +1.62%  "
$ws.Range("E50").Value = "'  +3.36%  "
$ws.Range("E51").Value = "'  -0.86%  "
